$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Update existing row 46: the XRP/USDT0000004 Buy order is now CANCEL, and
#    gets a Finalized date (I46) of 2017-05-08 06:49:21 (serial 42863.28427...)
# ---------------------------------------------------------------------------
$ws.Range("H46").Value = "CANCEL"
$ws.Range("I46").Value = 42863.284270833334

# ---------------------------------------------------------------------------
# 2) Append new row 47: a new Buy order for the same XRP/USDT0000004 batch,
#    created at the same timestamp as the finalized date above.
# ---------------------------------------------------------------------------

# Column A - Data/timestamp (date format with wrap text, same as A46/I46)
$ws.Range("A47").NumberFormat = "m/d/yy h:mm"
$ws.Range("A47").WrapText = $true
$ws.Range("A47").Value = 42863.284270833334

# Column B - Action (rich text "Buy" in green, matching the style used
# throughout the sheet for buy orders: 12 leading spaces + green "Buy").
$ws.Range("B47").Value = "            Buy"
$ws.Range("B47").Characters(13, 3).Font.Color = 5287936

# Column C - Currency
$ws.Range("C47").Value = "        XRP"

# Column D - Current value (USDT) - numeric-looking text, must stay text.
# Build it as a formula first (forces text type) then flatten to a static
# shared-string value, and enable wrap text (matches D46's style).
$ws.Range("D47").Formula = "=""             0.1926" + [char]10 + """"
$ws.Range("D47").Value = $ws.Range("D47").Value
$ws.Range("D47").WrapText = $true

# Column E - Transaction value
$ws.Range("E47").Value = "         0.185USDT"

# Column F - Transaction amount
$ws.Range("F47").Value = "         210 XRP"

# Column G - Transaction code (same batch as row 46)
$ws.Range("G47").Value = " XRP/USDT0000004"

# Column H - Status
$ws.Range("H47").Value = "IN PROGRESS"

# Column I - Finalized date left blank (order still in progress), but keep
# the same date-formatted style as the other Finalized-date cells.
$ws.Range("I47").NumberFormat = "m/d/yy h:mm"
$ws.Range("I47").WrapText = $true

# Column K - Profit(%) placeholder (same blank-ish value used elsewhere)
$ws.Range("K47").Value = "     "

Write-Output "Row 46 updated and row 47 appended"
